$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'maa://24702 (94.16), maa://25390 (95.82), maa://36681 (85.92)'
$ws.Range("AB2").Value = 'maa://21246 (91.32), maa://36684 (97.56), ***maa://22731 (6.67)'
$ws.Range("AF2").Value = 'maa://25251 (92.05), ***maa://21730 (16.92), ***maa://39501 (16.67), *maa://36675 (60.0)'
$ws.Range("T4").Value = 'maa://32509 (97.85), maa://27295 (82.46), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range("X4").Value = '**maa://32495 (47.27), ***maa://31785 (22.22), ***maa://36683 (28.26), *maa://43217 (80.0)'
$ws.Range("D6").Value = 'maa://42407 (93.1)'
$ws.Range("AF6").Value = '*maa://33152 (60.0), ***maa://22770 (27.27)'
$ws.Range("H7").Value = '*maa://22763 (67.86)'
$ws.Range("A8").Value = '更新日期：2024.11.14 13:18:07'
$ws.Range("P8").Value = 'maa://32931 (84.21), *maa://21916 (60.0), maa://23252 (92.42), **maa://22759 (45.45), maa://37496 (95.83)'
$ws.Range("T9").Value = '**maa://22866 (30.77), maa://26222 (97.62)'
$ws.Range("X9").Value = 'maa://26223 (97.27)'
$ws.Range("AF9").Value = 'maa://26206 (89.8), **maa://22865 (48.98)'
$ws.Range("P10").Value = 'maa://28977 (91.25), maa://36669 (86.21), *maa://23264 (61.82)'
$ws.Range("T10").Value = 'maa://27395 (95.73), maa://22755 (87.39), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range("X10").Value = 'maa://22301 (97.54), maa://22726 (100.0)'
$ws.Range("AF10").Value = '*maa://25021 (54.43), *maa://22733 (59.38), maa://22761 (100.0)'
$ws.Range("L11").Value = 'maa://21287 (88.04)'
$ws.Range("T11").Value = 'maa://22747 (93.2), maa://22501 (98.21)'
$ws.Range("X11").Value = 'maa://36713 (98.06)'
$ws.Range("AF11").Value = 'maa://31203 (95.65), ***maa://24394 (19.23)'
$ws.Range("X12").Value = 'maa://22753 (91.5), *maa://21485 (76.87), maa://37962 (86.36)'
$ws.Range("AB12").Value = 'maa://23669 (95.26), maa://36677 (93.33), maa://39872 (90.0)'
$ws.Range("H13").Value = '*maa://21248 (74.54), **maa://22728 (47.73)'
$ws.Range("X13").Value = '*maa://34957 (77.59), *maa://22768 (51.61)'
$ws.Range("AF13").Value = '**maa://22737 (30.37), maa://39883 (90.91), *maa://39885 (58.33)'
$ws.Range("L14").Value = 'maa://26245 (96.18), maa://21288 (96.21), maa://36682 (97.3), maa://39841 (93.55)'
$ws.Range("AB14").Value = 'maa://22764 (96.67)'
$ws.Range("H15").Value = 'maa://24304 (88.66), maa://21478 (91.18)'
$ws.Range("D16").Value = 'maa://21441 (96.21), maa://36679 (92.68), maa://37650 (96.77)'
$ws.Range("T16").Value = 'maa://22729 (95.27), *maa://28648 (68.42), maa://36674 (83.78)'
$ws.Range("D18").Value = 'maa://24570 (96.92)'
$ws.Range("H18").Value = 'maa://24421 (89.87)'
$ws.Range("L18").Value = 'maa://22466 (88.89), *maa://22732 (50.6)'
$ws.Range("X18").Value = 'maa://21917 (97.65), maa://22741 (83.33)'
$ws.Range("AF18").Value = '*maa://24313 (57.32), **maa://29784 (44.44)'
$ws.Range("AB19").Value = '*maa://30709 (62.28), *maa://36668 (54.17)'
$ws.Range("D20").Value = 'maa://21432 (90.51), maa://25198 (92.86), *maa://20795 (50.4), maa://36680 (96.43)'
$ws.Range("H20").Value = 'maa://22864 (88.41)'
$ws.Range("L20").Value = 'maa://41331 (84.93)'
$ws.Range("H21").Value = 'maa://24372 (96.7)'
$ws.Range("P21").Value = 'maa://24381 (86.67)'
$ws.Range("AB21").Value = '*maa://21443 (79.59), ***maa://23820 (29.82)'
$ws.Range("AF21").Value = 'maa://22524 (94.33), *maa://22432 (76.27)'
$ws.Range("H22").Value = 'maa://25236 (96.2), **maa://21678 (48.94), **maa://22735 (42.86)'
$ws.Range("X22").Value = 'maa://21282 (98.36), *maa://37649 (68.18)'
$ws.Range("AF22").Value = 'maa://29658 (92.86)'
$ws.Range("L23").Value = 'maa://39756 (93.07), maa://39875 (94.74)'
$ws.Range("AB23").Value = 'maa://29652 (97.44)'
$ws.Range("D24").Value = 'maa://24368 (80.23)'
$ws.Range("X24").Value = 'maa://29988 (86.3), maa://23504 (92.95), **maa://22892 (39.86), *maa://25141 (77.42), maa://36663 (80.65), ***maa://22815 (23.08)'
$ws.Range("AF24").Value = 'maa://22523 (85.42), *maa://36672 (79.59), maa://29910 (92.31), **maa://21440 (34.55)'
$ws.Range("D25").Value = 'maa://29753 (95.04)'
$ws.Range("H25").Value = '*maa://29063 (74.47), *maa://25311 (75.26), ***maa://22725 (4.84)'
$ws.Range("T25").Value = 'maa://20109 (92.17), maa://22545 (100.0), maa://42915 (100.0)'
$ws.Range("AB25").Value = 'maa://31215 (85.71), *maa://24516 (79.07), maa://26001 (87.27)'
$ws.Range("AF25").Value = 'maa://20108 (96.21), maa://24621 (96.55), maa://36676 (96.43), maa://22771 (85.71), maa://37772 (100.0)'
$ws.Range("AB26").Value = 'maa://42235 (91.38)'
$ws.Range("T28").Value = 'maa://23263 (94.79), *maa://29765 (60.27)'
$ws.Range("X28").Value = 'maa://39929 (88.93), ***maa://39723 (14.29), maa://41749 (82.76)'
$ws.Range("AF28").Value = 'maa://36660 (92.39), *maa://36701 (62.96)'
$ws.Range("L29").Value = 'maa://28432 (93.33), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range("AF29").Value = '*maa://24080 (69.25), ***maa://34960 (8.7), maa://42865 (86.96)'
$ws.Range("AB30").Value = 'maa://42979 (96.49)'
$ws.Range("T32").Value = 'maa://41108 (87.5), maa://42859 (92.86), maa://41238 (95.83)'
$ws.Range("T34").Value = 'maa://24526 (93.39)'
$ws.Range("L35").Value = 'maa://41296 (95.29)'
$ws.Range("T36").Value = 'maa://27613 (98.99)'
$ws.Range("P38").Value = '*maa://24383 (67.74)'
$ws.Range("AF38").Value = 'maa://36697 (85.26)'
$ws.Range("P40").Value = 'maa://23278 (95.92), maa://21386 (95.7), maa://36664 (90.0)'
$ws.Range("H43").Value = 'maa://22525 (92.25), maa://21284 (83.33)'
$ws.Range("H45").Value = 'maa://21229 (85.08), maa://30807 (95.24), *maa://22767 (57.89), ***maa://20796 (13.79), *maa://42459 (60.0)'
$ws.Range("T45").Value = '**maa://39364 (41.18)'
$ws.Range("H47").Value = 'maa://27410 (95.98), maa://29661 (97.76), maa://28038 (84.62)'
$ws.Range("H55").Value = 'maa://32532 (92.08)'
$ws.Range("H57").Value = 'maa://25176 (98.11)'
$ws.Range("H59").Value = 'maa://27746 (83.5), maa://31270 (95.41)'
$ws.Range("H60").Value = '*maa://40438 (54.29)'
